$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 70 (shifts existing rows 70-108 down to 71-109)
$ws.Rows.Item(70).Insert()

# Populate the newly inserted row 70 with this week's price report
$ws.Range("A70").Value = 2
$ws.Range("B70").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C70").Value = "Coquimbo"
$ws.Range("D70").Value = 45176
$ws.Range("E70").Value = 4
$ws.Range("F70").Value = 100112022
$ws.Range("G70").Value = "Arveja Verde"
$ws.Range("H70").Value = "Perfection"
$ws.Range("I70").Value = "Primera"
$ws.Range("J70").Value = 240
$ws.Range("K70").Value = 21000
$ws.Range("L70").Value = 23000
$ws.Range("M70").Value = 22000
$ws.Range("N70").Value = "$/malla 25 kilos"
$ws.Range("O70").Value = "Provincia de Limarí"
$ws.Range("P70").Value = 880
$ws.Range("Q70").Value = 25
$ws.Range("R70").Value = "Hortaliza"
